# This script updates the "想去人数" (want-to-go count) column (F) for a
# handful of rows on both the "展览" and "全部类型" worksheets, matching
# the re-scraped values from the upstream data source.

$wb = $excel.ActiveWorkbook

# Map of row number -> new F-column value (identical update on both sheets).
$updates = @{
    3  = 142
    4  = 1385
    5  = 1618
    6  = 356
    7  = 468
    9  = 195
    16 = 1793
    23 = 4333
    24 = 17
    28 = 47
    29 = 692
    31 = 348
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
